$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix the capitalisation of the signatory's surname:
#    "Sammie degerlien" -> "Sammie Degerlien"
#    This phrase occurs twice in the document (once near the top in a plain
#    run, once near the bottom in the bold signature block). Only the
#    second (bold) occurrence must be corrected, so we locate both matches
#    and only touch the one whose run is bold.
# ---------------------------------------------------------------------------

$searchText = "Sammie degerlien"

# First occurrence (plain / not bold) - search from the start of the document.
$firstRange = $d.Content
$firstFound = $firstRange.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Second occurrence (bold) - search from right after the first match.
$secondRange = $d.Range($firstRange.End, $d.Content.End)
$secondFound = $secondRange.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($secondFound -and $secondRange.Font.Bold) {
    # Change only the lower-case "d" (at offset 7 within "Sammie degerlien")
    # to an upper-case "D", leaving every other character/run untouched.
    $dStart = $secondRange.Start + 7
    $insertionPoint = $d.Range($dStart, $dStart)
    $insertionPoint.InsertBefore("D")

    # The original lower-case "d" has now shifted one position to the right.
    $oldLetter = $d.Range($dStart + 1, $dStart + 2)
    $oldLetter.Delete()
}

# ---------------------------------------------------------------------------
# 2. Reposition the logo picture anchored in the page header.
#    positionH (from margin): 4276725 EMU -> 3905250 EMU  (336.75pt -> 307.5pt)
#    positionV (from paragraph): 9525 EMU -> 19050 EMU      (0.75pt  -> 1.5pt)
# ---------------------------------------------------------------------------

$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)
$logo = $header.Shapes.Item(1)
$logo.Left = 307.5
$logo.Top = 1.5
